$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# These cells hold numeric-looking values that must remain TEXT (shared
# strings), matching the source file's original storage. A leading
# apostrophe forces Excel to keep the literal text instead of silently
# coercing it to a number; re-applying the "Normal" style afterwards keeps
# the cell's formatting/style index untouched (the apostrophe trick alone
# would tag the cell with a derived quote-prefixed style).
$ws.Range("B13").Value = "'9.37"
$ws.Range("B13").Style = "Normal"

$ws.Range("C13").Value = "'1.74"
$ws.Range("C13").Style = "Normal"

$ws.Range("D13").Value = "'11.11"
$ws.Range("D13").Style = "Normal"

$ws.Range("B15").Value = "'81.72"
$ws.Range("B15").Style = "Normal"

$ws.Range("C15").Value = "'15.17"
$ws.Range("C15").Style = "Normal"

$ws.Range("D15").Value = "'96.88"
$ws.Range("D15").Style = "Normal"
